$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.922.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.876.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7384"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.93%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.68"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.04%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3146"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07227"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.63"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.17%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08342"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7491"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.892.83"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.388"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.53%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.24"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.916.08"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.095"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "246.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007835"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.140.76"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.22%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.009"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.02%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1552"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.294"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.24"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.64"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.021"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.494"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.606"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.09%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.250"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05326"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.233"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.91%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7489"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9990"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01959"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.756"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4524"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.113.55"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.134"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.40"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8630"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.30"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.002"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.860"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.580"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.522"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.038.70"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.61%  "
